# Update the "Sheet1" book-price table with refreshed scrape data.
# (Sheet2 keeps the previous snapshot of the same data and is left untouched.)

$wb = $excel.ActiveWorkbook
$ws1 = $wb.ActiveSheet

# Row 2 - price refreshed
$ws1.Range("C2").Value = "₹570.13"

# Row 3 - new entry: Learning Robotic Process Automation / Alok Mani Tripathi
$ws1.Range("A3").Value = "Learning Robotic Process Automation: Create Software robots and automate business processes with the leading RPA tool – UiPath"
$ws1.Range("B3").Value = " Alok Mani Tripathi  "
$ws1.Range("C3").Value = "₹2,971"

# Row 4 - Robotic Process Automation Projects / Nandan Mullakara
$ws1.Range("A4").Value = "Robotic Process Automation Projects: Build real-world RPA solutions using UiPath and Automation Anywhere"
$ws1.Range("B4").Value = " Nandan Mullakara and Arun Kumar Asokan  "
$ws1.Range("C4").Value = "₹390.58"

# Row 5 - RPA for everyone / Vincenzo Marchica
$ws1.Range("A5").Value = "RPA for everyone: Robotic Process Automation, this famous unknown"
$ws1.Range("B5").Value = " Vincenzo Marchica "
$ws1.Range("C5").Value = "₹0"

# Row 6 - unchanged (Real World RPA Use Cases / Srikanth Merianda / ₹209)

# Row 7 - Introduction To RPA / Abhinav Sabharwal
$ws1.Range("A7").Value = "Introduction To RPA"
$ws1.Range("B7").Value = " Abhinav Sabharwal  "
$ws1.Range("C7").Value = "₹0"

# Row 8 - RPA Tutorial / Terrell Tromburg
$ws1.Range("A8").Value = "RPA Tutorial: A Guide To Learning Rpa For The Average Worker To Preserve Your Job"
$ws1.Range("B8").Value = " Terrell Tromburg "
$ws1.Range("C8").Value = "₹0"

# Row 9 - RPA (Japanese Edition) / Author not Available
$ws1.Range("A9").Value = "RPA (Japanese Edition)"
$ws1.Range("B9").Value = "Author not Available"
$ws1.Range("C9").Value = "₹0"

# Row 10 - RPA for Developers / Stijn Van Hijfte (no price listed)
$ws1.Range("A10").Value = "RPA for Developers : Automation Solutions Every Developer Should Know"
$ws1.Range("B10").Value = " Stijn Van Hijfte "
$ws1.Range("C10").ClearContents()

# Row 11 - Robotic Process Automation (RPA): Business Presentation / FlevyPro Library
$ws1.Range("A11").Value = "Robotic Process Automation (RPA): Business Presentation (FlevyPro Frameworks)"
$ws1.Range("B11").Value = " FlevyPro Library "
$ws1.Range("C11").Value = "₹0"
